# Rebuild the body as a list of "first test" write-ups, each consisting of a
# short text blurb followed (after a blank line) by a Google-Drive hyperlink.
#
# The original document had a single paragraph:
#   "Map Test 2.0 " + hyperlink(...bk05VXdLMEZqeHM)
# followed by a paragraph that only carries the _GoBack bookmark.
#
# We delete the old paragraph's content (merging it away so the bookmark
# paragraph becomes paragraph 1), then rebuild all of the new paragraphs in
# front of it by repeatedly inserting text/paragraph-breaks/hyperlinks at a
# tracked position. Finally we append one more empty paragraph after the
# bookmark paragraph (before the sectPr).

$d = $word.ActiveDocument

function Insert-TextPara($doc, $pos, $prevEnd, $text) {
    $str = $text + "`r"
    $r = $doc.Range($pos, $pos)
    $r.InsertBefore($str)
    $newEnd = $doc.Content.End
    $newPos = $pos + ($newEnd - $prevEnd)
    return $newPos, $newEnd
}

function Insert-EmptyPara($doc, $pos, $prevEnd) {
    $r = $doc.Range($pos, $pos)
    $r.InsertBefore("`r")
    $newEnd = $doc.Content.End
    $newPos = $pos + ($newEnd - $prevEnd)
    return $newPos, $newEnd
}

function Insert-LinkPara($doc, $pos, $prevEnd, $url) {
    # Reserve the paragraph mark for the hyperlink's own paragraph first so
    # Hyperlinks.Add has somewhere to stop instead of merging into whatever
    # paragraph happens to follow.
    $mark = $doc.Range($pos, $pos)
    $mark.InsertBefore("`r")
    $hr = $doc.Range($pos, $pos)
    $doc.Hyperlinks.Add($hr, $url, $null, $null, $url) | Out-Null
    $newEnd = $doc.Content.End
    $newPos = $pos + ($newEnd - $prevEnd)
    return $newPos, $newEnd
}

# 1. Blow away the old paragraph's content. This merges paragraph 1 into
#    paragraph 2 (the bookmark-only paragraph), which now becomes
#    paragraph 1, ready to have new content inserted in front of it.
$p1 = $d.Paragraphs.Item(1)
$delRange = $d.Range($p1.Range.Start, $p1.Range.End)
$delRange.Delete()

# 2. The ordered list of new content blocks.
$items = @(
    @{ Type = "text"; Text = "First combat test: " },
    @{ Type = "link"; Url = "https://drive.google.com/open?id=0B0dYxrDwUlTxZi0xMjRKakhSQUk" },
    @{ Type = "empty" },

    @{ Type = "text"; Text = "Equipment test:" },
    @{ Type = "link"; Url = "https://drive.google.com/open?id=0B0dYxrDwUlTxcHVZZTBUZ3A4Mlk" },
    @{ Type = "empty" },

    @{ Type = "text"; Text = "First map test (obsolete, algorithm changed)" },
    @{ Type = "link"; Url = "https://drive.google.com/open?id=0B0dYxrDwUlTxNVBvdGktMkhwdnc" },
    @{ Type = "empty" },

    @{ Type = "text"; Text = "Map Test 2.0, " },
    @{ Type = "text"; Text = "improved algorithm, faster generation, objects inside zones" },
    @{ Type = "link"; Url = "https://drive.google.com/open?id=0B0dYxrDwUlTxbk05VXdLMEZqeHM" },
    @{ Type = "empty" },

    @{ Type = "text"; Text = "Map Tesp 2.5, " },
    @{ Type = "text"; Text = "Allows movement between world and locations" },
    @{ Type = "link"; Url = "https://drive.google.com/open?id=0B0dYxrDwUlTxNzdmLWxYY25acms" },
    @{ Type = "empty" },

    @{ Type = "text"; Text = "Map Test 2.7" },
    @{ Type = "text"; Text = "Added transition between world and location and vice versa, movement between locations and world seperated" },
    @{ Type = "link"; Url = "https://drive.google.com/open?id=0B0dYxrDwUlTxckNUeW4zTGlfUEk" }
)

$pos = 0
$prevEnd = $d.Content.End

foreach ($it in $items) {
    if ($it.Type -eq "text") {
        $res = Insert-TextPara $d $pos $prevEnd $it.Text
    } elseif ($it.Type -eq "empty") {
        $res = Insert-EmptyPara $d $pos $prevEnd
    } else {
        $res = Insert-LinkPara $d $pos $prevEnd $it.Url
    }
    $pos = $res[0]
    $prevEnd = $res[1]
}

# 3. Append one more empty paragraph after the bookmark paragraph (which is
#    now the very last content paragraph before sectPr).
$d.Paragraphs.Add() | Out-Null

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
Write-Output ("Hyperlinks.Count=" + $d.Hyperlinks.Count)
